$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "Approved " status values from D2 and D3 (the shared string
# becomes unused and is dropped from sharedStrings.xml as a result).
$ws.Range("D2:D3").ClearContents()

# Select D3 to match the final selection left in the sheet view.
$ws.Range("D3").Select()
